$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 96 (shifts existing rows 96+ down by one,
# and Excel copies formatting from the row above, i.e. style index 2,
# matching the target workbook's s="2" styling for the new row).
$ws.Rows.Item(96).Insert()

# Populate the newly inserted row 96 with the new ontology term.
$ws.Range("A96").Value = "GMHO:0000173"
$ws.Range("B96").Value = "post-traumatic stress disorder severity"
$ws.Range("C96").Value = "A data item that is about the location on the dimension of post-traumatic stress disorder."
$ws.Range("D96").Value = "anxiety disorder severity"
$ws.Range("G96").Value = "PTSD severity"
$ws.Range("S96").Value = "Proposed"
$ws.Range("V96").Value = "PS"
